$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price/Volume cells store numeric-looking values as plain text
# (inline strings in the source workbook). Setting NumberFormat to "@"
# (Text) on each target cell before assigning its new value keeps Excel
# from auto-converting the string into a Number/Percentage.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.33%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "15.42%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.085"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.70%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07869"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.83%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.294"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.11%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.091"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.91%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.003"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.74%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9251"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.69%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1006"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.16%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.42%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08654"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.97%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03413"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09906"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.90%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001480"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.02%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04662"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.62%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005597"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.80%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.27%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.63%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3434"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.13%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1319"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.45%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.558"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.63%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.83%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.08%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004491"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.69%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.31%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002998"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-11.44%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.55%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04699"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.02%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007841"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.22%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1418"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.82%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002298"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.78%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009182"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.75%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006003"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.74%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.800"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "123.60%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002688"
